# Update forecast_summary_B0BH9DXY38.xlsx with corrected forecast output.
#
# 1. "Forecast Comparison" sheet:
#    - Insert a new "Week_Start_Date" column between "Week" and "ASIN".
#    - Normalize the Week labels ("W01" -> "W1", etc.).
#    - Correct a handful of MyForecast values.
#    - Store is_holiday_week as a boolean instead of a number.
# 2. "Summary" sheet:
#    - Refresh the derived forecast statistics to match the corrected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert the new Week_Start_Date column (becomes column B) ---------------
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)
# MyForecast values (column D after the insert) — corrected figures.
$myForecast = @(226,202,198,204,191,194,206,232,203,203,204,242,246,192,192,192)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Week label without the leading zero.
    $ws.Cells.Item($row, 1).Value = $weeks[$i]

    # Week_Start_Date as plain text (leading apostrophe keeps Excel from
    # auto-converting the ISO date string into a date serial number).
    $ws.Cells.Item($row, 2).Value = "'" + $weekStartDates[$i]

    # Corrected MyForecast figure.
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]

    # is_holiday_week now stored as a boolean.
    $ws.Cells.Item($row, 10).Value = $false
}

# --- Refresh the Summary sheet statistics ------------------------------------
# Leading apostrophes keep these as plain text (matching the existing
# "numbers/dates stored as strings" convention on this sheet) instead of
# letting Excel auto-coerce them into numeric or date values.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "'3326"
$summary.Range("B10").Value = "'1653"
$summary.Range("B11").Value = "'830"
$summary.Range("B12").Value = "'246"
$summary.Range("B13").Value = "'2025-03-30"
$summary.Range("B14").Value = "'191"
